$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 "Save" - copy style from neighboring header cell G1 so it
# reuses the same bold/border/centered header style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value2 = "Save"

# Fill H2:H49 with a 0/1 flag: 1 when the row's "sum" (column G) exceeds 8,
# else 0.
$lastRow = 49
for ($r = 2; $r -le $lastRow; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2
    if ($g -gt 8) {
        $ws.Cells.Item($r, 8).Value2 = 1
    } else {
        $ws.Cells.Item($r, 8).Value2 = 0
    }
}
